$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 1
    6  = 1
    7  = 0
    8  = 1
    9  = 0
    10 = 1
    11 = 4
    12 = 2
    13 = 1
    14 = 0
    15 = 0
    16 = 0
    17 = 1
    18 = 1
    19 = 3
    20 = 1
    21 = 3
    22 = 0
    23 = 0
    24 = 1
    25 = 1
    26 = 1
    27 = 0
    28 = 0
    29 = 0
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
